$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 additions: Ended time, What was Confirmed, What was Done
# (written in this order so shared-string table indices land the same way Excel produced them)
$ws.Range("D12").Value = "9:30PM"
$ws.Range("F12").Value = "Chatrooms UID are based off of the users in the room (it is created when the first message is sent), User ID are based on the first letter of the first name and the first letter of their last name followed by 5 digits "
$ws.Range("E12").Value = "Worked on Presentation slides and class diagram - Adding Structure of Presentation"

# Match formatting used by the rest of the table: D (time) matches column C/D's time style,
# E/F (notes) match column E's text style
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("F12").PasteSpecial(-4122)

# Widen column E to fit the new, longer "What was Done" text
# (target stored width is 83.83203125; the engine's ColumnWidth setter quantizes
# to 1/6-character steps, so this is the closest reachable setting)
$ws.Columns("E").ColumnWidth = 82.99869791666667
